# Applies the commit "calculation of new indicators":
#  - SCHEME_MEASURES: renumber MQMS0x -> MQME00x
#  - METADATA_ISSUES: remap rule codes to the new MQME0xx scheme
#  - METADATA_MEASURES: drop the old "Total number of columns" row and
#    renumber the remaining rows
#  - METADATA_METRICS: replace the IQME0x indicators with the new MQID0xx
#    set (7 renamed/updated rows + 4 brand-new rows)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# SCHEME_MEASURES  (MQMS01..MQMS05 -> MQME001..MQME005)
# ---------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsMeasures.Range("A2").Value = "MQME001"
$wsMeasures.Range("A3").Value = "MQME002"
$wsMeasures.Range("A4").Value = "MQME003"
$wsMeasures.Range("A5").Value = "MQME004"
$wsMeasures.Range("A6").Value = "MQME005"

# ---------------------------------------------------------------------
# METADATA_ISSUES  (rule-code remap, rows 2..490)
# ---------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")
$lastRow = 490
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsIssues.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -eq "MQME10") {
        $cell.Value = "MQME012"
    } elseif ($val -eq "MQME12") {
        $cell.Value = "MQME014"
    } elseif ($val -eq "MQME01") {
        $cell.Value = "MQME008"
    } elseif ($val -eq "MQME14") {
        $cell.Value = "MQME009"
    } elseif ($val -eq "MQME15") {
        $cell.Value = "MQME010"
    }
}

# ---------------------------------------------------------------------
# METADATA_MEASURES
#   old row2 "MQME00  / Total number of columns / 1042" is removed
#   old row3 "MQMEA1  / Total number of length-required columns / 423"
#            becomes new row2 "MQME006"
#   old row4 "MQMEA2  / Total number of NUMBER columns / 476"
#            becomes new row3 "MQME007"
# ---------------------------------------------------------------------
$wsMetaMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMetaMeasures.Range("A2").Value = "MQME006"
$wsMetaMeasures.Range("B2").Value = "Total number of length-required columns"
$wsMetaMeasures.Range("C2").Value = 423

$wsMetaMeasures.Range("A3").Value = "MQME007"
$wsMetaMeasures.Range("B3").Value = "Total number of NUMBER columns"
$wsMetaMeasures.Range("C3").Value = 476

$wsMetaMeasures.Range("A4").ClearContents()
$wsMetaMeasures.Range("B4").ClearContents()
$wsMetaMeasures.Range("C4").ClearContents()

# ---------------------------------------------------------------------
# METADATA_METRICS
#   replace IQME01..IQME07 with the new MQID001..MQID011 indicators
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")

$wsMetrics.Range("A2").Value = "MQID001"
$wsMetrics.Range("B2").Value = "Table names in singular"
$wsMetrics.Range("C2").Value = "98.89%"

$wsMetrics.Range("A3").Value = "MQID002"
$wsMetrics.Range("B3").Value = "Table with recommended name length"
$wsMetrics.Range("C3").Value = "100.00%"

$wsMetrics.Range("A4").Value = "MQID003"
$wsMetrics.Range("B4").Value = "Columns with correct prefixes"
$wsMetrics.Range("C4").Value = "91.17%"

$wsMetrics.Range("A5").Value = "MQID004"
$wsMetrics.Range("B5").Value = "Columns with recommended name size"
$wsMetrics.Range("C5").Value = "100.00%"

$wsMetrics.Range("A6").Value = "MQID005"
$wsMetrics.Range("B6").Value = "Columns with comments"
$wsMetrics.Range("C6").Value = "66.89%"

$wsMetrics.Range("A7").Value = "MQID006"
$wsMetrics.Range("B7").Value = "Table with standard PK prefixes"
$wsMetrics.Range("C7").Value = "64.79%"

$wsMetrics.Range("A8").Value = "MQID007"
$wsMetrics.Range("B8").Value = "Table with standard FK prefixes"
$wsMetrics.Range("C8").Value = "61.19%"

$wsMetrics.Range("A9").Value = "MQID008"
$wsMetrics.Range("B9").Value = "Table with standard UK prefixes"
$wsMetrics.Range("C9").Value = "100.00%"

$wsMetrics.Range("A10").Value = "MQID009"
$wsMetrics.Range("B10").Value = "NUMBER columns with valid scale"
$wsMetrics.Range("C10").Value = "100.00%"

$wsMetrics.Range("A11").Value = "MQID010"
$wsMetrics.Range("B11").Value = "Columns with valid num_distinct"
$wsMetrics.Range("C11").Value = "100.00%"

$wsMetrics.Range("A12").Value = "MQID011"
$wsMetrics.Range("B12").Value = "Columns with valid num_nulls"
$wsMetrics.Range("C12").Value = "100.00%"
